# Adds/updates 2024-03-14 daily crime-index data across all sheets.
# Generated from the authoritative diff of output/cta-index-crime-ytd.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("E2").Value = 11
$ws.Range("J2").Value = 30
$ws.Range("D3").Value = 29
$ws.Range("F3").Value = 20
$ws.Range("J3").Value = 32
$ws.Range("F9").Value = 112
$ws.Range("G9").Value = 108
$ws.Range("I9").Value = 98
$ws.Range("J9").Value = 82
$ws.Range("K9").Value = 79
$ws.Range("B10").Value = 186
$ws.Range("C10").Value = 225
$ws.Range("D10").Value = 352
$ws.Range("E10").Value = 344
$ws.Range("F10").Value = 495
$ws.Range("G10").Value = 390
$ws.Range("I10").Value = 146
$ws.Range("J10").Value = 116
$ws.Range("K10").Value = 148
$ws.Range("B11").Value = 300
$ws.Range("C11").Value = 344
$ws.Range("D11").Value = 478
$ws.Range("E11").Value = 467
$ws.Range("F11").Value = 644
$ws.Range("G11").Value = 545
$ws.Range("I11").Value = 296
$ws.Range("J11").Value = 268
$ws.Range("K11").Value = 287

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 5

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("C5").Value = 12
$ws.Range("C6").Value = 19

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("F7").Value = 6
$ws.Range("F8").Value = 11

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 7
$ws.Range("F3").Value = 2
$ws.Range("J3").Value = 3
$ws.Range("B7").Value = 20
$ws.Range("C7").Value = 35
$ws.Range("D7").Value = 104
$ws.Range("E7").Value = 93
$ws.Range("I7").Value = 42
$ws.Range("K7").Value = 24
$ws.Range("B8").Value = 30
$ws.Range("C8").Value = 47
$ws.Range("D8").Value = 120
$ws.Range("E8").Value = 110
$ws.Range("F8").Value = 161
$ws.Range("I8").Value = 70
$ws.Range("J8").Value = 44
$ws.Range("K8").Value = 41

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K5").Value = 4
$ws.Range("C6").Value = 4
$ws.Range("G6").Value = 8
$ws.Range("C7").Value = 4
$ws.Range("G7").Value = 12
$ws.Range("K7").Value = 8

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K5").Value = 6
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("C5").Value = 4
$ws.Range("F5").Value = 12
$ws.Range("C6").Value = 9
$ws.Range("F6").Value = 18

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("G2").Value = 4
$ws.Range("G7").Value = 17
$ws.Range("K20").Value = 5
$ws.Range("K26").Value = 4
$ws.Range("C31").Value = 19
$ws.Range("F35").Value = 11
$ws.Range("F40").Value = 3
$ws.Range("I40").Value = 2
$ws.Range("E42").Value = 3
$ws.Range("G44").Value = 7
$ws.Range("K49").Value = 9
$ws.Range("B52").Value = 30
$ws.Range("C52").Value = 47
$ws.Range("D52").Value = 120
$ws.Range("E52").Value = 110
$ws.Range("F52").Value = 161
$ws.Range("I52").Value = 70
$ws.Range("J52").Value = 44
$ws.Range("K52").Value = 41
$ws.Range("G55").Value = 3
$ws.Range("E60").Value = 10
$ws.Range("J61").Value = 4
$ws.Range("C64").Value = 9
$ws.Range("F64").Value = 18
$ws.Range("G66").Value = 3
$ws.Range("F67").Value = 10
$ws.Range("C69").Value = 4
$ws.Range("G69").Value = 12
$ws.Range("K69").Value = 8
$ws.Range("E71").Value = 3
$ws.Range("C73").Value = 4
$ws.Range("D73").Value = 12
$ws.Range("D75").Value = 11
$ws.Range("K75").Value = 12
$ws.Range("I76").Value = 12
$ws.Range("G78").Value = 9
$ws.Range("F81").Value = 3
$ws.Range("G84").Value = 4
$ws.Range("J84").Value = 3
$ws.Range("K84").Value = 5
$ws.Range("E85").Value = 7
$ws.Range("C93").Value = 5
$ws.Range("J93").Value = 2
$ws.Range("F94").Value = 4
$ws.Range("B97").Value = 300
$ws.Range("C97").Value = 344
$ws.Range("D97").Value = 478
$ws.Range("E97").Value = 467
$ws.Range("F97").Value = 644
$ws.Range("G97").Value = 545
$ws.Range("I97").Value = 296
$ws.Range("J97").Value = 268
$ws.Range("K97").Value = 287

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("G5").Value = 9
$ws.Range("G6").Value = 9

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("F5").Value = 2
$ws.Range("H6").Value = 1
$ws.Range("F7").Value = 3
$ws.Range("H7").Value = 2

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("E6").Value = 5
$ws.Range("E7").Value = 7

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("G5").Value = 6
$ws.Range("G6").Value = 7

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("D6").Value = 10
$ws.Range("K6").Value = 6
$ws.Range("D7").Value = 11
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = 3

$ws = $wb.Worksheets.Item('River North')
$ws.Range("D3").Value = 1
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 12

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("G4").Value = 2
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 4
$ws.Range("G6").Value = 4
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 5

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 3

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K5").Value = 2
$ws.Range("K6").Value = 4
$ws.Range("J6").Value = 2
$ws.Range("J7").Value = 4

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 4

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I6").Value = 4
$ws.Range("I8").Value = 12

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J4").Value = 1
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 5
$ws.Range("J6").Value = 2

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("F3").Value = 1
$ws.Range("F6").Value = 4

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("G5").Value = 14
$ws.Range("G6").Value = 17

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 3

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("F6").Value = 9
$ws.Range("F7").Value = 10

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = 3

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("C2").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 3

